$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: dmqm_seminar entry update
$ws.Range("D36").Value = "Graph-based semi-supervised learning"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/313"

# Row 37: dsba_seminar entry update
$ws.Range("D37").Value = "[Paper Review] SOM-DST  : Efficient Dialogue State Tracking by Selectively Overwriting Memory"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1449&mod=document&pageid=1"

# Row 41: cloudinsight entry update
$ws.Range("D41").Value = "오픈 포맷 데이터 성능 향상 전략 – Part 1 캐싱"
$ws.Range("E41").Value = "http://cloudinsight.net/data/%ec%98%a4%ed%94%88-%ed%8f%ac%eb%a7%b7-%eb%8d%b0%ec%9d%b4%ed%84%b0-%ec%84%b1%eb%8a%a5-%ed%96%a5%ec%83%81-%ec%a0%84%eb%9e%b5-part-1-%ec%ba%90%ec%8b%b1/"

# Row 51: bskyvsion entry update
$ws.Range("D51").Value = "[css] 링크 밑줄 없애기"
$ws.Range("E51").Value = "https://bskyvision.com/1145"
